# Update TPM-derived metric columns (E..T) on the active sheet to reflect
# the newly recomputed values described by the commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 5.516293333333333
$ws.Cells.Item(2, 8).Value = 16.54888
$ws.Cells.Item(2, 9).Value = 0.02672612669241994
$ws.Cells.Item(2, 10).Value = 0.02672612669241994
$ws.Cells.Item(2, 13).Value = 0.004819666666666667
$ws.Cells.Item(2, 14).Value = 0.014459
$ws.Cells.Item(2, 15).Value = 0.04945987179224048
$ws.Cells.Item(2, 16).Value = 0.04945987179224049
$ws.Cells.Item(2, 17).Value = 0.02658669510222222
$ws.Cells.Item(2, 18).Value = 0.23928025592
$ws.Cells.Item(2, 19).Value = 0.001321870799710267
$ws.Cells.Item(2, 20).Value = 0.001321870799710267
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 5.516293333333333
$ws.Cells.Item(3, 8).Value = 16.54888
$ws.Cells.Item(3, 9).Value = 0.02672612669241994
$ws.Cells.Item(3, 10).Value = 0.02672612669241994
$ws.Cells.Item(3, 15).Value = 0.5155094445470654
$ws.Cells.Item(3, 16).Value = 0.5155094445470654
$ws.Cells.Item(3, 17).Value = 0.2771073180711111
$ws.Cells.Item(3, 18).Value = 2.49396586264
$ws.Cells.Item(3, 19).Value = 0.0137775707261039
$ws.Cells.Item(3, 20).Value = 0.0137775707261039
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 5.516293333333333
$ws.Cells.Item(4, 8).Value = 16.54888
$ws.Cells.Item(4, 9).Value = 0.02672612669241994
$ws.Cells.Item(4, 10).Value = 0.02672612669241994
$ws.Cells.Item(4, 13).Value = 0.04239200000000001
$ws.Cells.Item(4, 14).Value = 0.127176
$ws.Cells.Item(4, 15).Value = 0.4350306836606942
$ws.Cells.Item(4, 16).Value = 0.4350306836606942
$ws.Cells.Item(4, 17).Value = 0.2338467069866667
$ws.Cells.Item(4, 18).Value = 2.10462036288
$ws.Cells.Item(4, 19).Value = 0.01162668516660578
$ws.Cells.Item(4, 20).Value = 0.01162668516660578
$ws.Cells.Item(5, 9).Value = 0.02157503014446814
$ws.Cells.Item(5, 10).Value = 0.02157503014446814
$ws.Cells.Item(5, 13).Value = 0.004819666666666667
$ws.Cells.Item(5, 14).Value = 0.014459
$ws.Cells.Item(5, 15).Value = 0.04945987179224048
$ws.Cells.Item(5, 16).Value = 0.04945987179224049
$ws.Cells.Item(5, 17).Value = 0.02146247209233333
$ws.Cells.Item(5, 18).Value = 0.193162248831
$ws.Cells.Item(5, 19).Value = 0.001067098224859118
$ws.Cells.Item(5, 20).Value = 0.001067098224859118
$ws.Cells.Item(6, 9).Value = 0.02157503014446814
$ws.Cells.Item(6, 10).Value = 0.02157503014446814
$ws.Cells.Item(6, 15).Value = 0.5155094445470654
$ws.Cells.Item(6, 16).Value = 0.5155094445470654
$ws.Cells.Item(6, 19).Value = 0.01112213180586096
$ws.Cells.Item(6, 20).Value = 0.01112213180586096
$ws.Cells.Item(7, 9).Value = 0.02157503014446814
$ws.Cells.Item(7, 10).Value = 0.02157503014446814
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.04239200000000001
$ws.Cells.Item(7, 14).Value = 0.127176
$ws.Cells.Item(7, 15).Value = 0.4350306836606942
$ws.Cells.Item(7, 16).Value = 0.4350306836606942
$ws.Cells.Item(7, 17).Value = 0.188775942376
$ws.Cells.Item(7, 18).Value = 1.698983481384
$ws.Cells.Item(7, 19).Value = 0.009385800113748061
$ws.Cells.Item(7, 20).Value = 0.009385800113748061
$ws.Cells.Item(8, 7).Value = 196.4313813333333
$ws.Cells.Item(8, 8).Value = 589.294144
$ws.Cells.Item(8, 9).Value = 0.9516988431631119
$ws.Cells.Item(8, 10).Value = 0.9516988431631119
$ws.Cells.Item(8, 13).Value = 0.004819666666666667
$ws.Cells.Item(8, 14).Value = 0.014459
$ws.Cells.Item(8, 15).Value = 0.04945987179224048
$ws.Cells.Item(8, 16).Value = 0.04945987179224049
$ws.Cells.Item(8, 17).Value = 0.9467337808995555
$ws.Cells.Item(8, 18).Value = 8.520604028095999
$ws.Cells.Item(8, 19).Value = 0.0470709027676711
$ws.Cells.Item(8, 20).Value = 0.04707090276767111
$ws.Cells.Item(9, 7).Value = 196.4313813333333
$ws.Cells.Item(9, 8).Value = 589.294144
$ws.Cells.Item(9, 9).Value = 0.9516988431631119
$ws.Cells.Item(9, 10).Value = 0.9516988431631119
$ws.Cells.Item(9, 15).Value = 0.5155094445470654
$ws.Cells.Item(9, 16).Value = 0.5155094445470654
$ws.Cells.Item(9, 17).Value = 9.867599487025776
$ws.Cells.Item(9, 18).Value = 88.808395383232
$ws.Cells.Item(9, 19).Value = 0.4906097420151005
$ws.Cells.Item(9, 20).Value = 0.4906097420151005
$ws.Cells.Item(10, 7).Value = 196.4313813333333
$ws.Cells.Item(10, 8).Value = 589.294144
$ws.Cells.Item(10, 9).Value = 0.9516988431631119
$ws.Cells.Item(10, 10).Value = 0.9516988431631119
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.04239200000000001
$ws.Cells.Item(10, 14).Value = 0.127176
$ws.Cells.Item(10, 15).Value = 0.4350306836606942
$ws.Cells.Item(10, 16).Value = 0.4350306836606942
$ws.Cells.Item(10, 17).Value = 8.327119117482667
$ws.Cells.Item(10, 18).Value = 74.944072057344
$ws.Cells.Item(10, 19).Value = 0.4140181983803404
$ws.Cells.Item(10, 20).Value = 0.4140181983803404
